# Refactor "Triggers" sheet name to "Processors"
$wb = $excel.ActiveWorkbook

$sheet = $wb.Worksheets.Item("Triggers")
$sheet.Name = "Processors"

# Select cell C23 on the renamed (formerly "Triggers") sheet, as that
# sheet was active in the source edit and its selection moved.
$sheet.Activate()
$sheet.Range("C23").Select()
